# Refresh cryptocurrency price / 1h-volume figures (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.004.27"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "'2.043.09"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'246.23"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "'57.28"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'0.0770"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "'15.59"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "'0.887"
$ws.Range("E13").Value = "  +11.79%  "
$ws.Range("D14").Value = "'2.346.79"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "'5.68"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'2.050.72"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  +11.17%  "
$ws.Range("D18").Value = "'36.983.34"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'74.20"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'5.42"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "'235.62"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'170.00"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.48"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D28").Value = "'19.92"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "'5.41"
$ws.Range("E29").Value = "  +15.43%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").Value = "'4.71"
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("D33").Value = "'0.0613"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "'0.0874"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.84"
$ws.Range("E36").Value = "  +5.65%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.26"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "'5.08"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "'0.0987"
$ws.Range("E41").Value = "  -7.25%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'1.14"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").Value = "'97.25"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "'1.294.06"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "'2.36"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").Value = "'2.86"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").Value = "'3.73"
$ws.Range("E49").Value = "  +6.80%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'6.78"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "'2.230.20"
$ws.Range("E51").Value = "  -0.36%  "
